$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The serializer now emits a "user" field right after "id", so the
# "user" column (previously the last data column, P) needs to move to
# column C, pushing dsc..get_user_application_detail one column to the
# right (C..O -> D..P). Columns Q, R, S (in_*_approved_by) stay put.

# 1) Insert a new blank column at C. This shifts the old C..S to D..T,
#    so the old "user" column (P) now lives at Q.
$ws.Columns("C").Insert()

# 2) Move the "user" header/value (now at Q1:Q2) into the freshly
#    inserted C1:C2, then remove the now-empty Q1:Q2, which shifts the
#    trailing in_district/in_state/in_central columns back from R,S,T
#    (T was never populated) to Q,R,S.
$ws.Range("Q1:Q2").Copy($ws.Range("C1:C2"))
$ws.Range("Q1:Q2").Delete()

# 3) The updated_at value (now at K2) reflects a newer timestamp.
$ws.Range("K2").Value = "2022-11-30T13:26:43.523007+05:45"
